# national_parks.xlsx - "map photo filepath update"
#
# Sequence of events reconstructed from the target diff:
#   1) While Sheet1 was scrolled/selected at H39, the user duplicated the
#      sheet ("Move or Copy" -> Create a copy), producing "Sheet1 (2)" as
#      a pristine snapshot of the *unedited* data (still pointing at the
#      old local /Users/alexreed/... image file paths).
#   2) Back on Sheet1, the local file-path strings in column H (one per
#      park photo) were replaced with the equivalent public
#      raw.githubusercontent.com URLs, turned into real hyperlinks, ending
#      with the selection on the last edited cell, H62.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1) snapshot the sheet before touching anything -------------------
$ws.Range("H39").Select() | Out-Null
$ws.Copy([Type]::Missing, $ws) | Out-Null

# Hop back onto the original sheet (the copy just activated "Sheet1 (2)").
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate() | Out-Null

# --- 2) swap each local image path for its public URL + hyperlink -----
$photoLinks = @(
    @{ Cell = "H4";  Url = "https://raw.githubusercontent.com/reedalexandria/reedalexandria.github.io/main/image/parks/arches.jpeg" },
    @{ Cell = "H9";  Url = "https://raw.githubusercontent.com/reedalexandria/reedalexandria.github.io/main/image/parks/bryce.JPG" },
    @{ Cell = "H24"; Url = "https://raw.githubusercontent.com/reedalexandria/reedalexandria.github.io/main/image/parks/glacier.jpeg" },
    @{ Cell = "H26"; Url = "https://raw.githubusercontent.com/reedalexandria/reedalexandria.github.io/main/image/parks/teton.jpeg" },
    @{ Cell = "H32"; Url = "https://raw.githubusercontent.com/reedalexandria/reedalexandria.github.io/main/image/parks/volcano.jpg" },
    @{ Cell = "H39"; Url = "https://raw.githubusercontent.com/reedalexandria/reedalexandria.github.io/main/image/parks/kings.jpeg" },
    @{ Cell = "H45"; Url = "https://raw.githubusercontent.com/reedalexandria/reedalexandria.github.io/main/image/parks/rainier.JPG" },
    @{ Cell = "H53"; Url = "https://raw.githubusercontent.com/reedalexandria/reedalexandria.github.io/main/image/parks/sequoia.jpg" },
    @{ Cell = "H61"; Url = "https://raw.githubusercontent.com/reedalexandria/reedalexandria.github.io/main/image/parks/yosemite.JPG" },
    @{ Cell = "H62"; Url = "https://raw.githubusercontent.com/reedalexandria/reedalexandria.github.io/main/image/parks/zion.JPG" }
)

foreach ($link in $photoLinks) {
    $ws.Hyperlinks.Add($ws.Range($link.Cell), $link.Url, [Type]::Missing, [Type]::Missing, $link.Url) | Out-Null
}

# Leave the cursor where the last edit happened.
$ws.Range("H62").Select() | Out-Null
